$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 422; existing rows 422:497 shift down to 423:498.
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with the new price-record data.
$ws.Range("A422").Value = 10
$ws.Range("B422").Value = "Vega Modelo de Temuco"
$ws.Range("C422").Value = "La Araucanía"
$ws.Range("D422").Value = 44505
$ws.Range("E422").Value = 9
$ws.Range("F422").Value = 100112021
$ws.Range("G422").Value = "Ají"
$ws.Range("H422").Value = "Inferno"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 120
$ws.Range("K422").Value = 34000
$ws.Range("L422").Value = 35000
$ws.Range("M422").Value = 34542
$ws.Range("N422").Value = "$/caja 15 kilos"
$ws.Range("O422").Value = "Región de Arica y Parinacota"
$ws.Range("P422").Value = 2303
$ws.Range("Q422").Value = 15
$ws.Range("R422").Value = "Hortaliza"
